# Add a new "Note" column (F) to the BOM sheet, with a manufacturing /
# moisture-sensitivity-level note attached to the USB-C receptacle row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 - bold, matching the rest of row 1's header style.
$ws.Range("F1").Value = "Note"
$ws.Range("F1").Font.Bold = $true

# Note text for row 7 (GCT USB4105-GF-A USB-C receptacle, J1).
$note = @"
The products sealed in moisture barrier bags (MBB) should be stored in a non-condensing atmospheric
environment of < 40 °C and 90%RH. The module is rated at the moisture sensitivity level (MSL) of 3.
After unpacking, the module must be soldered within 168 hours with the factory conditions 25±5 °C and
60%RH. If the above conditions are not met, the module needs to be baked.
"@

$ws.Range("F7").Value = $note
$ws.Range("F7").WrapText = $true
$ws.Range("F7").Font.Bold = $false

# Grow row 7 so the wrapped note text is fully visible.
$ws.Rows.Item(7).RowHeight = 46.25

# Mirror the author's updated viewport/selection from the commit.
$excel.ActiveWindow.ScrollColumn = 3 | Out-Null
$ws.Range("E60").Select() | Out-Null
